$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.704.27'
$ws.Range('E2').Value = '  +0.58%  '
$ws.Range('D3').Value = '3.226.84'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '609.55'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.31%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '159.04'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.20%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.225.68'
$ws.Range('E8').Value = '  +1.12%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.552'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('E11').Value = '  -4.75%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.504'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.58%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000271'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '38.90'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').Value = '3.755.72'
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D16').Value = '66.714.40'
$ws.Range('E16').Value = '  +0.57%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.37'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('D18').Value = '3.225.34'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '508.00'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.18'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.32%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.735'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.02'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '84.94'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '3.01'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.13'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('E30').Value = '  +34.68%  '
$ws.Range('B31').Value = 'Stacks'
$ws.Range('C31').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.94'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.02'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.95%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '28.17'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('E35').Value = '  -3.86%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.47'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '502.81'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -1.34%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '55.53'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('D39').Value = '0.0₃0774'
$ws.Range('E39').Value = '  +15.45%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.133'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.47%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.09'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +7.41%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0421'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.74%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.72'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.298'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.98%  '
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('D46').Value = '2.906.11'
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '28.22'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.63%  '
$ws.Range('E48').Value = '  +3.10%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '122.24'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.03%  '
